$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column B costs (row 2-31) from original currency to new currency
# by dividing by the exchange rate 29.5, and drop the custom cell style
# (border/Arial font) those cells had, reverting them to the default style.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $orig = $cell.Value2
    $cell.ClearFormats()
    $cell.Value2 = $orig / 29.5
}

# Update the active selection to D32 (single cell), as recorded in the sheet view.
$ws.Range("D32").Select()
